# Rename the *img sheets to img* (img prefix instead of img suffix).
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

# The renamed "imge" sheet (formerly "eimg", the last tab) becomes the
# active/selected sheet, replacing "holiday" as the active tab.
$wb.Worksheets.Item("imge").Activate()
